# Swap the full row content between pairs of rows: (3,4), (5,6), (17,18).
# The underlying edit is a pure data re-ordering: every field that was in
# row N ends up in row M and vice versa, columns A..AY. Columns Y and AA
# hold a literal text date ("2026-01-24") that is identical before/after
# in every touched row, so they are intentionally left untouched (this
# also sidesteps Excel's automatic text->date coercion on write).
# Column I occasionally holds a numeric-looking text value ("1"); it is
# copied separately with the destination cell pre-formatted as Text so
# Excel does not silently coerce it into a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowRange($r1, $r2, $colStart, $colEnd) {
    $addr1 = "$colStart$r1" + ":" + "$colEnd$r1"
    $addr2 = "$colStart$r2" + ":" + "$colEnd$r2"
    $rng1 = $ws.Range($addr1)
    $rng2 = $ws.Range($addr2)

    $val1 = $rng1.Value2
    $val2 = $rng2.Value2

    $rng1.Value = $val2
    $rng2.Value = $val1
}

function Swap-Cell($r1, $r2, $col) {
    $cell1 = $ws.Range("$col$r1")
    $cell2 = $ws.Range("$col$r2")

    $v1 = $cell1.Value2
    $v2 = $cell2.Value2

    # Force Text format so numeric-looking strings ("1", etc.) round-trip
    # as text instead of being coerced into numbers.
    $cell1.NumberFormat = "@"
    $cell2.NumberFormat = "@"

    if ($v2 -eq $false -or $v2 -eq $null) {
        $cell1.Value = ""
    } else {
        $cell1.Value = [string]$v2
    }

    if ($v1 -eq $false -or $v1 -eq $null) {
        $cell2.Value = ""
    } else {
        $cell2.Value = [string]$v1
    }
}

function Swap-FullRow($r1, $r2) {
    # A through H
    Swap-RowRange $r1 $r2 "A" "H"
    # Column I is handled specially (text-safe).
    Swap-Cell $r1 $r2 "I"
    # J through X
    Swap-RowRange $r1 $r2 "J" "X"
    # Columns Y (and Z) and AA are left untouched on purpose: identical
    # before/after, and re-assigning them risks date auto-coercion.
    # AB through AY
    Swap-RowRange $r1 $r2 "AB" "AY"
}

Swap-FullRow 3 4
Swap-FullRow 5 6
Swap-FullRow 17 18

Write-Host "Row swaps complete"
